# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap Kenia/Gabon order (Kenia now comes before Gabon) and refresh their stats ---
# Row 115 previously held Gabon's data; it now becomes Kenia with updated figures.
$ws.Range("A115").Value = "Kenia"
$ws.Range("B115").Value = 672
$ws.Range("C115").Value = 23
$ws.Range("D115").Value = 239
$ws.Range("E115").Value = 401
$ws.Range("F115").Value = 1
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 32

# Row 116 previously held Kenia's data; it now becomes Gabon, reusing the old Gabon figures.
$ws.Range("A116").Value = "Gabon"
$ws.Range("B116").Value = 661
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 110
$ws.Range("E116").Value = 543
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 8

# --- Update statistics for other countries/provinces ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1348315
$ws.Range("C4").Value = 1006
$ws.Range("E4").Value = 1030179
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 80056

# Row 17: India
$ws.Range("B17").Value = 64139
$ws.Range("C17").Value = 1331
$ws.Range("E17").Value = 42667
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 2114

# Row 45: Serbia
$ws.Range("B45").Value = 10114
$ws.Range("C45").Value = 82
$ws.Range("D45").Value = 3006
$ws.Range("E45").Value = 6893
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 215

# Row 139: Zambia
$ws.Range("B139").Value = 267
$ws.Range("C139").Value = 15
$ws.Range("D139").Value = 117
$ws.Range("E139").Value = 143
